$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.079.79"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "1.877.00"
$ws.Range("E3").Value = "  -2.33%  "
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.14%  "
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5039"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.69%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3958"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.32%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08217"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.04"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.093"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.56"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.13%  "
$ws.Range("D13").Value = "1.877.20"
$ws.Range("E13").Value = "  -1.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.296"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.201"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.78"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001086"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06469"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.02%  "
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("D22").Value = "30.079.54"
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("E23").Value = "  -3.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.154"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.60%  "
$ws.Range("D26").Value = "2.100.85"
$ws.Range("E26").Value = "  -1.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.244"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.59"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.074"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("E32").Value = "  -2.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.921"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.694"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02428"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.275"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06352"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2132"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.172"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.504"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.221"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.41%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6288"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.21%  "
$ws.Range("E44").Value = "  +0.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.24"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5905"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.087"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.632"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "122.28"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.99%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.208"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "77.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.05%  "
